$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.097257733345032
$ws.Range("B1").Value = 1.103580117225647
$ws.Range("C1").Value = 1.066441655158997
$ws.Range("D1").Value = 1.28650951385498
$ws.Range("E1").Value = 1.212653756141663
